$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 391, pushing the existing row 391 (and everything
# below it) down by one. This is what grows the used range from A1:R477 to
# A1:R478.
$ws.Rows.Item(391).Insert()

# Populate the newly inserted row 391 with the new record.
$ws.Cells.Item(391, 1).Value2 = 7
$ws.Cells.Item(391, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(391, 3).Value2 = "Ñuble"
$ws.Cells.Item(391, 4).Value2 = 45135
$ws.Cells.Item(391, 5).Value2 = 16
$ws.Cells.Item(391, 6).Value2 = 100112006
$ws.Cells.Item(391, 7).Value2 = "Repollo"
$ws.Cells.Item(391, 8).Value2 = "Crespo record"
$ws.Cells.Item(391, 9).Value2 = "Primera"
$ws.Cells.Item(391, 10).Value2 = 150
$ws.Cells.Item(391, 11).Value2 = 1000
$ws.Cells.Item(391, 12).Value2 = 1000
$ws.Cells.Item(391, 13).Value2 = 1000
$ws.Cells.Item(391, 14).Value2 = "$/unidad"
$ws.Cells.Item(391, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(391, 16).Value2 = 1000
$ws.Cells.Item(391, 17).Value2 = 1
$ws.Cells.Item(391, 18).Value2 = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D
# (the row-insert already copies the neighbouring format, this just makes
# sure it's explicit/idempotent).
$ws.Cells.Item(391, 4).NumberFormat = $ws.Cells.Item(392, 4).NumberFormat
